$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Heading: merge "IT " + "Work" (spell-checked) into a single clean run
#    "IT Work" with no w:proofErr wrapping.
# ---------------------------------------------------------------------------
$headingPara = $d.Paragraphs(1)
$headingRange = $headingPara.Range
$headingRange.InsertParagraphBefore()
$newHeadingPara = $d.Paragraphs(1)
$newHeadingPara.Range.Text = "IT Work"
$newHeadingPara.Style = $d.Styles("Heading1")
$oldHeadingPara = $d.Paragraphs(2)
$oldHeadingPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2) "...what's broken. " -> "...what's broken" + closing curly quote + moved
#    _GoBack bookmark + ". "
# ---------------------------------------------------------------------------
$quote = [char]0x201D
$findRng = $d.Content
$findRng.Find.Execute(("what" + [char]0x2019 + "s broken"), $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$brokenEnd = $findRng.End
$insertPoint = $d.Range($brokenEnd, $brokenEnd)
$insertPoint.Text = $quote
$bookmarkPos = $brokenEnd + 1
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# ---------------------------------------------------------------------------
# 3) "engaging their existing" -> "engaging the clients existing"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("engaging their existing", $true, $false, $false, $false, $false, $true, 1, $false, "engaging the clients existing", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) The old _GoBack bookmark location (near "solutioning mode") is cleaned
#    up automatically because _GoBack is a single, unique bookmark per
#    document - re-adding it above already moved it off of that spot.
# ---------------------------------------------------------------------------
